# Fruta / hortaliza, semanal
# Insert a new weekly price record for Pomelo (Start Ruby, Primera) at row 238,
# pushing the existing rows 238-332 down to 239-333.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 238 - this shifts rows 238..332 down
# to 239..333 (and copies formatting from the row above, matching row style).
$ws.Rows.Item(238).Insert()

# Populate the newly inserted row 238 with the new weekly record.
$ws.Cells.Item(238, 1).Value  = 10
$ws.Cells.Item(238, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(238, 3).Value  = "La Araucanía"
$ws.Cells.Item(238, 4).Value  = 44875
$ws.Cells.Item(238, 5).Value  = 9
$ws.Cells.Item(238, 6).Value  = "Fruta"
$ws.Cells.Item(238, 7).Value  = 100102
$ws.Cells.Item(238, 8).Value  = "Cítricos"
$ws.Cells.Item(238, 9).Value  = 100102006
$ws.Cells.Item(238, 10).Value = "Pomelo"
$ws.Cells.Item(238, 11).Value = "Start Ruby"
$ws.Cells.Item(238, 12).Value = "Primera"
$ws.Cells.Item(238, 13).Value = 240
$ws.Cells.Item(238, 14).Value = 12000
$ws.Cells.Item(238, 15).Value = 13000
$ws.Cells.Item(238, 16).Value = 12583
$ws.Cells.Item(238, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(238, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(238, 19).Value = 839
$ws.Cells.Item(238, 20).Value = 15
